# SampleInventory.xlsx — add new "Nathan's Eye Bags (Unwashed)" row to the
# Merchandise sheet, tweak the Zach's Croc unit price, and bring the few
# stray "no-decimal currency" cells in line with the rest of the UnitPrice
# column (2-decimal currency formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Merchandise")

# --- Normalize UnitPrice formatting -------------------------------------
# C7, C10 and C12 were still using the old "$#,##0" (no decimals) number
# format while every other cell in the UnitPrice column uses "$#,##0.00".
# Re-apply the same 2-decimal currency format so they match.
$ws.Range("C7").NumberFormat = """$""#,##0.00"
$ws.Range("C10").NumberFormat = """$""#,##0.00"
$ws.Range("C12").NumberFormat = """$""#,##0.00"

# --- Update Zach's Croc stock count -------------------------------------
$ws.Range("C12").Value = 1

# --- Add the new inventory row ------------------------------------------
$ws.Range("A14").Value = "Nathan's Eye Bags (Unwashed)"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 2
$ws.Range("C14").NumberFormat = """$""#,##0.00"

# --- Match the saved selection -------------------------------------------
$ws.Range("I13").Select()
